# Updates Price (D) and Volume(1h) (E) columns on the cryptos sheet.
# Values are written as text (NumberFormat "@") and the cell style is
# reset to "Normal" afterwards so Excel does not coerce numeric-looking
# strings (e.g. "216.70") into actual numbers and no stray style index
# is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "26.099.50"
Set-TextValue $ws.Range("E2") "  -1.53%  "
Set-TextValue $ws.Range("D3") "1.667.39"
Set-TextValue $ws.Range("E3") "  -0.73%  "
Set-TextValue $ws.Range("E4") "  +0.05%  "
Set-TextValue $ws.Range("D5") "216.70"
Set-TextValue $ws.Range("E5") "  +0.17%  "
Set-TextValue $ws.Range("D6") "0.5110"
Set-TextValue $ws.Range("E6") "  +2.94%  "
Set-TextValue $ws.Range("E7") "  +0.09%  "
Set-TextValue $ws.Range("D8") "0.2632"
Set-TextValue $ws.Range("E8") "  +1.35%  "
Set-TextValue $ws.Range("D9") "0.06415"
Set-TextValue $ws.Range("E9") "  +3.75%  "
Set-TextValue $ws.Range("D10") "21.72"
Set-TextValue $ws.Range("E10") "  -0.22%  "
Set-TextValue $ws.Range("E11") "  +2.06%  "
Set-TextValue $ws.Range("D12") "1.665.90"
Set-TextValue $ws.Range("E12") "  -0.47%  "
Set-TextValue $ws.Range("D13") "4.510"
Set-TextValue $ws.Range("E13") "  +1.84%  "
Set-TextValue $ws.Range("E14") "  +1.38%  "
Set-TextValue $ws.Range("D15") "0.000008571"
Set-TextValue $ws.Range("E15") "  +5.04%  "
Set-TextValue $ws.Range("D16") "64.38"
Set-TextValue $ws.Range("E16") "  +0.13%  "
Set-TextValue $ws.Range("D17") "26.160.98"
Set-TextValue $ws.Range("E17") "  -1.27%  "
Set-TextValue $ws.Range("E18") "  -1.12%  "
Set-TextValue $ws.Range("D19") "1.005"
Set-TextValue $ws.Range("E19") "  -0.02%  "
Set-TextValue $ws.Range("E20") "  +0.35%  "
Set-TextValue $ws.Range("D21") "188.99"
Set-TextValue $ws.Range("E21") "  +2.77%  "
Set-TextValue $ws.Range("D22") "6.206"
Set-TextValue $ws.Range("E22") "  +0.55%  "
Set-TextValue $ws.Range("E23") "  +0.06%  "
Set-TextValue $ws.Range("D24") "145.87"
Set-TextValue $ws.Range("E24") "  +1.07%  "
Set-TextValue $ws.Range("D25") "7.632"
Set-TextValue $ws.Range("E25") "  +2.50%  "
Set-TextValue $ws.Range("D26") "0.1192"
Set-TextValue $ws.Range("E26") "  +5.66%  "
Set-TextValue $ws.Range("D27") "15.61"
Set-TextValue $ws.Range("E27") "  +1.42%  "
Set-TextValue $ws.Range("D28") "0.06428"
Set-TextValue $ws.Range("E28") "  +13.21%  "
Set-TextValue $ws.Range("D29") "1.307"
Set-TextValue $ws.Range("E29") "  +0.48%  "
Set-TextValue $ws.Range("D30") "1.319"
Set-TextValue $ws.Range("E30") "  -0.14%  "
Set-TextValue $ws.Range("D31") "3.525"
Set-TextValue $ws.Range("E31") "  +1.59%  "
Set-TextValue $ws.Range("D32") "3.514"
Set-TextValue $ws.Range("E32") "  +1.66%  "
Set-TextValue $ws.Range("D33") "1.637"
Set-TextValue $ws.Range("E33") "  +0.28%  "
Set-TextValue $ws.Range("D34") "1.019"
Set-TextValue $ws.Range("E34") "  +1.49%  "
Set-TextValue $ws.Range("E35") "  +3.13%  "
Set-TextValue $ws.Range("D36") "2.366"
Set-TextValue $ws.Range("E36") "  -0.04%  "
Set-TextValue $ws.Range("E37") "  +2.15%  "
Set-TextValue $ws.Range("D38") "6.198"
Set-TextValue $ws.Range("E38") "  +5.56%  "
Set-TextValue $ws.Range("D39") "0.01614"
Set-TextValue $ws.Range("E39") "  +1.89%  "
Set-TextValue $ws.Range("D40") "1.076.18"
Set-TextValue $ws.Range("E40") "  +0.79%  "
Set-TextValue $ws.Range("D41") "0.8610"
Set-TextValue $ws.Range("E41") "  +1.26%  "
Set-TextValue $ws.Range("D42") "1.010"
Set-TextValue $ws.Range("E42") "  +0.82%  "
Set-TextValue $ws.Range("D43") "100.76"
Set-TextValue $ws.Range("E43") "  +2.66%  "
Set-TextValue $ws.Range("D44") "1.815.19"
Set-TextValue $ws.Range("E44") "  -1.08%  "
Set-TextValue $ws.Range("D45") "0.00000000111"
Set-TextValue $ws.Range("E45") "  +7.44%  "
Set-TextValue $ws.Range("D46") "56.19"
Set-TextValue $ws.Range("E46") "  +0.29%  "
Set-TextValue $ws.Range("E47") "  +0.34%  "
Set-TextValue $ws.Range("D48") "8.059"
Set-TextValue $ws.Range("D49") "0.05207"
Set-TextValue $ws.Range("E49") "  +0.39%  "
Set-TextValue $ws.Range("E50") "  -0.38%  "
Set-TextValue $ws.Range("D51") "5.947"
Set-TextValue $ws.Range("E51") "  +5.83%  "
